$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-528). The commit bumped that date by one day for the entire
# column, i.e. 46081 -> 46082.
$rng = $ws.Range("C2:C528")
$rng.Value2 = 46082
